$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "categories_title"
$ws.Range("B16").Value = "CATEGORIES"
$ws.Range("A17").Value = "score"
$ws.Range("B17").Value = "Score:"
$ws.Range("A18").Value = "rank"
$ws.Range("B18").Value = "Rank:"
$ws.Range("C18").ClearContents()
$ws.Range("A19").Value = "property"
$ws.Range("B19").Value = "Property"
$ws.Range("C19").ClearContents()
$ws.Range("A20").Value = "polygon"
$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").ClearContents()
$ws.Range("A21").Value = "triangle"
$ws.Range("B21").Value = "Triangle"
$ws.Range("C21").Value = 1.5
$ws.Range("A22").Value = "quadrilateral"
$ws.Range("B22").Value = "Quadrilateral"
$ws.Range("C22").Value = 1.5
$ws.Range("A23").Value = "pentagon"
$ws.Range("B23").Value = "Pentagon"
$ws.Range("C23").Value = 1.5
$ws.Range("A24").Value = "hexagon"
$ws.Range("B24").Value = "Hexagon"
$ws.Range("C24").Value = 1.5
$ws.Range("A25").Value = "octagon"
$ws.Range("B25").Value = "Octagon"
$ws.Range("C25").Value = 1.5
$ws.Range("A26").Value = "decagon"
$ws.Range("B26").Value = "Decagon"
$ws.Range("C26").Value = 1.5
$ws.Range("A27").Value = "triangle_right"
$ws.Range("B27").Value = "Right Triangle"
$ws.Range("C27").Value = 2
$ws.Range("A28").Value = "triangle_equilateral"
$ws.Range("B28").Value = "Equilateral Triangle"
$ws.Range("C28").Value = 2
$ws.Range("A29").Value = "triangle_isosceles"
$ws.Range("B29").Value = "Isosceles Triangle"
$ws.Range("C29").Value = 2
$ws.Range("A30").Value = "triangle_scalene"
$ws.Range("B30").Value = "Scalene Triangle"
$ws.Range("C30").Value = 2
$ws.Range("A31").Value = "triangle_acute"
$ws.Range("B31").Value = "Acute Triangle"
$ws.Range("C31").Value = 2
$ws.Range("A32").Value = "triangle_obtuse"
$ws.Range("B32").Value = "Obtuse Triangle"
$ws.Range("C32").Value = 2
$ws.Range("A33").Value = "quad_rectangle"
$ws.Range("B33").Value = "Rectangle"
$ws.Range("C33").Value = 2
$ws.Range("A34").Value = "quad_square"
$ws.Range("B34").Value = "Square"
$ws.Range("C34").Value = 2
$ws.Range("A35").Value = "quad_rhombus"
$ws.Range("B35").Value = "Rhombus"
$ws.Range("C35").Value = 2
$ws.Range("A36").Value = "quad_parallelogram"
$ws.Range("B36").Value = "Parallelogram"
$ws.Range("C36").Value = 3
$ws.Range("A37").Value = "quad_trapezoid"
$ws.Range("B37").Value = "Trapezoid"
$ws.Range("C37").Value = 2.5
$ws.Range("A38").Value = "quad_trapezium"
$ws.Range("B38").Value = "Trapezium"
$ws.Range("C38").Value = 2.5
$ws.Range("A39").Value = "quad_kite"
$ws.Range("B39").Value = "Kite"
$ws.Range("C39").Value = 1
$ws.Range("A40").Value = "prop_side_3"
$ws.Range("B40").Value = "3 sides."
$ws.Range("C40").Value = 2
$ws.Range("A41").Value = "prop_side_4"
$ws.Range("B41").Value = "4 sides."
$ws.Range("C41").Value = 2
$ws.Range("A42").Value = "prop_side_5"
$ws.Range("B42").Value = "5 sides."
$ws.Range("C42").Value = 2
$ws.Range("A43").Value = "prop_side_6"
$ws.Range("B43").Value = "6 sides."
$ws.Range("C43").Value = 2
$ws.Range("A44").Value = "prop_side_8"
$ws.Range("B44").Value = "8 sides."
$ws.Range("C44").Value = 2
$ws.Range("A45").Value = "prop_side_10"
$ws.Range("B45").Value = "10 sides."
$ws.Range("C45").Value = 2
$ws.Range("A46").Value = "prop_90_degree"
$ws.Range("B46").Value = "Has a 90° angle."
$ws.Range("C46").Value = 3
$ws.Range("A47").Value = "prop_sides_equal_all"
$ws.Range("B47").Value = "All sides are equal."
$ws.Range("C47").Value = 3
$ws.Range("A48").Value = "prop_60_degree_all"
$ws.Range("B48").Value = "All angles equal 60°."
$ws.Range("C48").Value = 3
$ws.Range("A49").Value = "prop_sides_equal_two"
$ws.Range("B49").Value = "Two equal sides."
$ws.Range("C49").Value = 3
$ws.Range("A50").Value = "prop_angles_equal_two"
$ws.Range("B50").Value = "Two equal angles."
$ws.Range("C50").Value = 3
$ws.Range("A51").Value = "prop_sides_no_equal"
$ws.Range("B51").Value = "No sides are equal."
$ws.Range("C51").Value = 3
$ws.Range("A52").Value = "prop_angles_less_90_all"
$ws.Range("B52").Value = "All angles less than 90°."
$ws.Range("C52").Value = 3
$ws.Range("A53").Value = "prop_angle_greater_90"
$ws.Range("B53").Value = "Has an angle greater than 90°."
$ws.Range("C53").Value = 3
$ws.Range("A54").Value = "prop_90_degree_all"
$ws.Range("B54").Value = "All angles equal 90°."
$ws.Range("C54").Value = 3
$ws.Range("A55").Value = "prop_opposite_sides_parallel"
$ws.Range("B55").Value = "Opposite sides are parallel."
$ws.Range("C55").Value = 3
$ws.Range("A56").Value = "prop_opposite_sides_equal"
$ws.Range("B56").Value = "Opposite sides are equal."
$ws.Range("C56").Value = 3
$ws.Range("A57").Value = "prop_opposite_angles_equal"
$ws.Range("B57").Value = "Opposite angles are equal."
$ws.Range("C57").Value = 3
$ws.Range("A58").Value = "prop_opposite_pair_parallel"
$ws.Range("B58").Value = "Two sides are parallel."
$ws.Range("C58").Value = 3
$ws.Range("A59").Value = "prop_sides_no_parallel"
$ws.Range("B59").Value = "No sides are parallel."
$ws.Range("C59").Value = 3
$ws.Range("A60").Value = "prop_sides_pair_equal"
$ws.Range("B60").Value = "Two pairs of sides are equal."
$ws.Range("C60").Value = 3
$ws.Range("A61").Value = "shape_categories"
$ws.Range("B61").Value = "Shape Categories"
$ws.Range("A62").Value = "proceed_instruct"
$ws.Range("B62").Value = "Press this button to proceed."
$ws.Range("C62").Value = 3

# Rows 59-62 are new (beyond the previous used range A1:D58), so they don't
# automatically inherit the wrap-text style used by every other row in
# column B. Re-apply it explicitly.
$ws.Range("B59:B62").WrapText = $true

$ws.Range("B17").Select()
